$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '71.036.99'
$ws.Range('E2').Value = '  +6.24%  '

# Row 3
$ws.Range('D3').Value = '3.660.28'
$ws.Range('E3').Value = '  +6.30%  '

# Row 4
$cell = $ws.Range('D4')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '596.27'
$cell.ClearFormats()
$ws.Range('E5').Value = '  +2.80%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '194.67'
$cell.ClearFormats()
$ws.Range('E6').Value = '  +3.45%  '

# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.648'
$cell.ClearFormats()
$ws.Range('E7').Value = '  +2.75%  '

# Row 8
$ws.Range('D8').Value = '3.653.71'
$ws.Range('E8').Value = '  +6.25%  '

# Row 9
$ws.Range('E9').Value = '  -0.03%  '

# Row 10
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.185'
$cell.ClearFormats()
$ws.Range('E10').Value = '  +8.15%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.675'
$cell.ClearFormats()
$ws.Range('E11').Value = '  +4.80%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '58.45'
$cell.ClearFormats()
$ws.Range('E12').Value = '  +3.09%  '

# Row 13
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '0.0000296'
$cell.ClearFormats()
$ws.Range('E13').Value = '  +6.83%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '9.97'
$cell.ClearFormats()
$ws.Range('E14').Value = '  +6.09%  '

# Row 15
$ws.Range('D15').Value = '4.242.59'
$ws.Range('E15').Value = '  +6.24%  '

# Row 16
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '20.13'
$cell.ClearFormats()
$ws.Range('E16').Value = '  +7.59%  '

# Row 17
$ws.Range('D17').Value = '3.652.93'
$ws.Range('E17').Value = '  +5.55%  '

# Row 18
$ws.Range('D18').Value = '70.949.08'
$ws.Range('E18').Value = '  +6.05%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '12.85'
$cell.ClearFormats()
$ws.Range('E19').Value = '  +6.52%  '

# Row 20
$ws.Range('E20').Value = '  +2.77%  '

# Row 21
$ws.Range('E21').Value = '  +4.65%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '490.44'
$cell.ClearFormats()
$ws.Range('E22').Value = '  +1.79%  '

# Row 23
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '19.14'
$cell.ClearFormats()
$ws.Range('E23').Value = '  +13.61%  '

# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '5.30'
$cell.ClearFormats()
$ws.Range('E24').Value = '  +0.19%  '

# Row 25
$ws.Range('E25').Value = '  +4.26%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '91.59'
$cell.ClearFormats()
$ws.Range('E26').Value = '  +2.65%  '

# Row 27
$ws.Range('E27').Value = '  +7.05%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '11.53'
$cell.ClearFormats()
$ws.Range('E28').Value = '  +5.04%  '

# Row 29
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '9.64'
$cell.ClearFormats()
$ws.Range('E29').Value = '  +6.63%  '

# Row 30
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '32.93'
$cell.ClearFormats()
$ws.Range('E30').Value = '  +5.48%  '

# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '7.79'
$cell.ClearFormats()
$ws.Range('E31').Value = '  +6.35%  '

# Row 32
$ws.Range('E32').Value = '  +10.05%  '

# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '629.92'
$cell.ClearFormats()
$ws.Range('E33').Value = '  +5.52%  '

# Row 34
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '12.30'
$cell.ClearFormats()
$ws.Range('E34').Value = '  +4.65%  '

# Row 35
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '66.40'
$cell.ClearFormats()
$ws.Range('E35').Value = '  +4.13%  '

# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '40.31'
$cell.ClearFormats()
$ws.Range('E36').Value = '  +9.76%  '

# Row 37
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0831'
$ws.Range('E37').Value = '  +10.17%  '

# Row 38
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '0.415'
$cell.ClearFormats()
$ws.Range('E38').Value = '  +7.30%  '

# Row 39
$ws.Range('E39').Value = '  -0.01%  '

# Row 40
$ws.Range('E40').Value = '  +0.11%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '3.60'
$cell.ClearFormats()
$ws.Range('E41').Value = '  +1.78%  '

# Row 42
$ws.Range('D42').Value = '3.311.90'
$ws.Range('E42').Value = '  +2.33%  '

# Row 43
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '3.17'
$cell.ClearFormats()
$ws.Range('E43').Value = '  +9.54%  '

# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '2.84'
$cell.ClearFormats()
$ws.Range('E44').Value = '  +12.15%  '

# Row 45
$ws.Range('E45').Value = '  +6.75%  '

# Row 46
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '3.06'
$cell.ClearFormats()
$ws.Range('E46').Value = '  +7.61%  '

# Row 47
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '9.49'
$cell.ClearFormats()
$ws.Range('E47').Value = '  +10.06%  '

# Row 48
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.139'
$cell.ClearFormats()
$ws.Range('E48').Value = '  +3.84%  '

# Row 49
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '3.31'
$cell.ClearFormats()
$ws.Range('E49').Value = '  +3.08%  '

# Row 50
$ws.Range('E50').Value = '  -2.26%  '

# Row 51
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.ClearFormats()
$ws.Range('E51').Value = '  -0.15%  '

Write-Host "Applied cryptos update"